$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 2027
$ws.Range("I47").Value = 2027
$ws.Range("K47").Value = 2027
$ws.Range("M47").Value = -1055
$ws.Range("H64").Value = 4284.8
$ws.Range("J64").Value = 3960
$ws.Range("L64").Value = 3960
$ws.Range("N64").Value = -4456
$ws.Range("H67").Value = 4284.8
$ws.Range("J67").Value = 3960
$ws.Range("L67").Value = 3960
$ws.Range("N67").Value = -5676
$ws.Range("H74").Value = 9350
$ws.Range("J74").Value = 7904
$ws.Range("L74").Value = 7904
$ws.Range("N74").Value = -9776
$ws.Range("H77").Value = 9350
$ws.Range("J77").Value = 7904
$ws.Range("L77").Value = 39520
$ws.Range("N77").Value = -48880
$ws.Range("H107").Value = 1875.8823
$ws.Range("I107").Value = 1689.4546
$ws.Range("J107").Value = 2217.6667
$ws.Range("K107").Value = 1689.4546
$ws.Range("L107").Value = 2217.6667
$ws.Range("M107").Value = 230.5454
$ws.Range("N107").Value = -6057.6667
$ws.Range("H112").Value = 2177.7222
$ws.Range("J112").Value = 2177.7222
$ws.Range("L112").Value = 6533.1666
$ws.Range("N112").Value = -8749.1666
$ws.Range("H132").Value = 9529551
$ws.Range("I132").Value = 11907939
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 35723817
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -35721287
$ws.Range("N132").Value = -53060
$ws.Range("H137").Value = 1066.1613
$ws.Range("I137").Value = 1030.6923
$ws.Range("J137").Value = 1250.6
$ws.Range("K137").Value = 3092.0769
$ws.Range("L137").Value = 3751.8
$ws.Range("M137").Value = -542.0769
$ws.Range("N137").Value = -8851.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11025.7
$ws.Range("I2").Value = 1030.5
$ws.Range("K2").Value = 1030.5
$ws.Range("M2").Value = -917.5
$ws.Range("H61").Value = 1696.3334
$ws.Range("I61").Value = 1453.1
$ws.Range("J61").Value = 2182.8
$ws.Range("K61").Value = 1453.1
$ws.Range("L61").Value = 2182.8
$ws.Range("M61").Value = -1241.1
$ws.Range("N61").Value = -2606.8
$ws.Range("H74").Value = 1318.5714
$ws.Range("I74").Value = 979
$ws.Range("J74").Value = 1771.3334
$ws.Range("K74").Value = 979
$ws.Range("L74").Value = 1771.3334
$ws.Range("M74").Value = -105
$ws.Range("N74").Value = -3519.3334
$ws.Range("H77").Value = 1318.5714
$ws.Range("I77").Value = 979
$ws.Range("J77").Value = 1771.3334
$ws.Range("K77").Value = 4895
$ws.Range("L77").Value = 8856.666999999999
$ws.Range("M77").Value = -527
$ws.Range("N77").Value = -17592.667
$ws.Range("H88").Value = 2943.4285
$ws.Range("J88").Value = 2943.4285
$ws.Range("L88").Value = 2943.4285
$ws.Range("N88").Value = -3755.4285
$ws.Range("H91").Value = 2943.4285
$ws.Range("J91").Value = 2943.4285
$ws.Range("L91").Value = 2943.4285
$ws.Range("N91").Value = -5751.4285
$ws.Range("H116").Value = 11025.7
$ws.Range("I116").Value = 1030.5
$ws.Range("K116").Value = 1030.5
$ws.Range("M116").Value = 1263.5
$ws.Range("H132").Value = 3049.8235
$ws.Range("I132").Value = 2969.5833
$ws.Range("J132").Value = 3242.4
$ws.Range("K132").Value = 8908.749899999999
$ws.Range("L132").Value = 9727.200000000001
$ws.Range("M132").Value = -6378.749899999999
$ws.Range("N132").Value = -14787.2
$ws.Range("H136").Value = 1696.3334
$ws.Range("I136").Value = 1453.1
$ws.Range("J136").Value = 2182.8
$ws.Range("K136").Value = 4359.299999999999
$ws.Range("L136").Value = 6548.400000000001
$ws.Range("M136").Value = -1809.299999999999
$ws.Range("N136").Value = -11648.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11025.7
$ws.Range("I3").Value = 1030.5
$ws.Range("K3").Value = 1030.5
$ws.Range("M3").Value = -916.5
$ws.Range("H20").Value = 2459.8
$ws.Range("J20").Value = 1597.8
$ws.Range("L20").Value = 1597.8
$ws.Range("N20").Value = -2091.8
$ws.Range("H86").Value = 4534.278
$ws.Range("I86").Value = 5459.1665
$ws.Range("J86").Value = 2684.5
$ws.Range("K86").Value = 5459.1665
$ws.Range("L86").Value = 2684.5
$ws.Range("M86").Value = -4336.1665
$ws.Range("N86").Value = -4930.5
$ws.Range("H89").Value = 4534.278
$ws.Range("I89").Value = 5459.1665
$ws.Range("J89").Value = 2684.5
$ws.Range("K89").Value = 27295.8325
$ws.Range("L89").Value = 13422.5
$ws.Range("M89").Value = -21679.8325
$ws.Range("N89").Value = -24654.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1149.5834
$ws.Range("I58").Value = 1171
$ws.Range("J58").Value = 914
$ws.Range("K58").Value = 1171
$ws.Range("L58").Value = 914
$ws.Range("M58").Value = -968
$ws.Range("N58").Value = -1320
$ws.Range("H62").Value = 6899451.5
$ws.Range("I62").Value = 3003.7036
$ws.Range("K62").Value = 3003.7036
$ws.Range("M62").Value = -2379.7036
$ws.Range("H65").Value = 6899451.5
$ws.Range("I65").Value = 3003.7036
$ws.Range("K65").Value = 15018.518
$ws.Range("M65").Value = -11898.518
$ws.Range("H115").Value = 50249.25
$ws.Range("J115").Value = 50249.25
$ws.Range("L115").Value = 50249.25
$ws.Range("N115").Value = -52599.25
$ws.Range("H134").Value = 25642840
$ws.Range("I134").Value = 41668924
$ws.Range("J134").Value = 1102.8
$ws.Range("K134").Value = 125006772
$ws.Range("L134").Value = 3308.4
$ws.Range("M134").Value = -125004237
$ws.Range("N134").Value = -8378.4
$ws.Range("H136").Value = 1149.5834
$ws.Range("I136").Value = 1171
$ws.Range("J136").Value = 914
$ws.Range("K136").Value = 3513
$ws.Range("L136").Value = 2742
$ws.Range("M136").Value = -963
$ws.Range("N136").Value = -7842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 175
$ws.Range("I2").Value = 66.666664
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 399.999984
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -286.999984
$ws.Range("N2").Value = -3226
$ws.Range("H14").Value = 669.2
$ws.Range("I14").Value = 669.2
$ws.Range("K14").Value = 2007.6
$ws.Range("M14").Value = -1834.6
$ws.Range("H59").Value = 3549.9
$ws.Range("J59").Value = 4125
$ws.Range("L59").Value = 12375
$ws.Range("N59").Value = -13455
$ws.Range("H61").Value = 172.88889
$ws.Range("J61").Value = 237.25
$ws.Range("L61").Value = 711.75
$ws.Range("N61").Value = -1141.75
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 1185
$ws.Range("I107").Value = 698
$ws.Range("J107").Value = 1266.1666
$ws.Range("K107").Value = 2094
$ws.Range("L107").Value = 3798.4998
$ws.Range("M107").Value = -174
$ws.Range("N107").Value = -7638.4998
$ws.Range("H121").Value = 290.42856
$ws.Range("I121").Value = 290.42856
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 871.28568
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 438.71432
$ws.Range("N121").ClearContents()
$ws.Range("H131").Value = 28572886
$ws.Range("I131").Value = 111111390
$ws.Range("J131").Value = 1862.5769
$ws.Range("K131").Value = 333334170
$ws.Range("L131").Value = 5587.7307
$ws.Range("M131").Value = -333329130
$ws.Range("N131").Value = -15667.7307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18005802
$ws.Range("I70").Value = 16671782
$ws.Range("J70").Value = 20006830
$ws.Range("K70").Value = 16671782
$ws.Range("L70").Value = 20006830
$ws.Range("M70").Value = -16671512
$ws.Range("N70").Value = -20007370
$ws.Range("H73").Value = 18005802
$ws.Range("I73").Value = 16671782
$ws.Range("J73").Value = 20006830
$ws.Range("K73").Value = 16671782
$ws.Range("L73").Value = 20006830
$ws.Range("M73").Value = -16670846
$ws.Range("N73").Value = -20008702
$ws.Range("H80").Value = 3280
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 4760
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 4760
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -6756
$ws.Range("H83").Value = 3280
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 4760
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 23800
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -33784
$ws.Range("H122").Value = 2780
$ws.Range("I122").Value = 2750
$ws.Range("K122").Value = 8250
$ws.Range("M122").Value = -5800
$ws.Range("H126").Value = 2382.6
$ws.Range("I126").Value = 2078
$ws.Range("J126").Value = 2585.6667
$ws.Range("K126").Value = 6234
$ws.Range("L126").Value = 7757.000100000001
$ws.Range("M126").Value = -3764
$ws.Range("N126").Value = -12697.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7919.4375
$ws.Range("I136").Value = 10401.091
$ws.Range("J136").Value = 2459.8
$ws.Range("K136").Value = 31203.273
$ws.Range("L136").Value = 7379.400000000001
$ws.Range("M136").Value = -28653.273
$ws.Range("N136").Value = -12479.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3000.4
$ws.Range("I132").Value = 2300.25
$ws.Range("K132").Value = 6900.75
$ws.Range("M132").Value = -4370.75
$ws.Range("H133").Value = 29797.5
$ws.Range("J133").Value = 29797.5
$ws.Range("L133").Value = 29797.5
$ws.Range("N133").Value = -39917.5
$ws.Range("H136").Value = 717.9286
$ws.Range("I136").Value = 486.1
$ws.Range("K136").Value = 1458.3
$ws.Range("M136").Value = 1091.7
